# Apply the workbook edits described by the commit:
#  - add pet population and wages as regressors
#  - fix data page (Used Methods) tab-selection state
#  - fix summary page (Key_Assumptions) with correct assumptions / observations

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Key_Assumptions")
$ws2 = $wb.Worksheets.Item("Used Methods")

# --- Key_Assumptions sheet (sheet1) -----------------------------------------
# Row 2 (B2): now holds the regressors assumption text.
$ws1.Range("B2").Value = "External variables (wages, CPI, energy) used as regressors; correlations checked."

# Row 3 / Row 4: new Assumption texts (order matters so the new shared
# strings get appended in the same order as the target workbook).
$ws1.Range("B3").Value = "SARIMA was assumed to be the best model fore forecasting due to seasonality and upward trend."
$ws1.Range("B4").Value = 'Formulas for "Claims Cost" and "Premium" were assumed along with few variables.'

# Row 5: re-uses the existing "No sudden policy..." assumption text.
$ws1.Range("B5").Value = "No sudden policy or regulatory shocks assumed."

# Rows 6-8: keep their existing Observation text (type label updated below).
$ws1.Range("B6").Value = "Forecasts based on historical trends; extreme events may not be captured."
$ws1.Range("B7").Value = "Prophet model deviated from other forecasts, highlighting uncertainty in trend capture."
$ws1.Range("B8").Value = "Data before 2018 considered but may be less representative; trend checked with rolling windows."

# Row 9: replaced with a new Observation about wages.
$ws1.Range("B9").Value = "Wages has the highest impact; impact is positive; high significance."

# New rows 10 and 11: additional observations about pet population and
# the regressors not improving performance (text only for now).
$ws1.Range("B10").Value = "Pet population has little to no impact."
$ws1.Range("B11").Value = "Incoporating regressors did not improve model performance."

# Re-label the "Limitation" rows (and the two new rows) as "Observation"
# last, so the new "Observation" shared string is appended after all of
# the other new strings above.
$ws1.Range("A6").Value = "Observation"
$ws1.Range("A7").Value = "Observation"
$ws1.Range("A8").Value = "Observation"
$ws1.Range("A9").Value = "Observation"
$ws1.Range("A10").Value = "Observation"
$ws1.Range("A11").Value = "Observation"

# Copy the existing wrap-text/vertical-center style from A9 onto the new
# A10:A11 cells so they match the rest of the "Type" column formatting.
$ws1.Range("A9").Copy()
$ws1.Range("A10:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column B to fit the longer observation/assumption text.
$ws1.Columns.Item(2).ColumnWidth = 92.6166666667

# --- Used Methods sheet (sheet2) --------------------------------------------
# Content is unchanged; only the active-tab/selection state changes below.

# --- Active sheet / selection state -----------------------------------------
# Key_Assumptions becomes the active (visible) tab with B19 selected, and
# Used Methods keeps its own selection but is no longer the active tab.
$ws1.Activate()
$ws1.Range("B19").Select()
